# Refresh crypto price/volume snapshot (and a few re-ranked rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''59.916.10'
$ws.Range('E2').Value = '''  -1.02%  '

# Row 3
$ws.Range('D3').Value = '''2.663.79'
$ws.Range('E3').Value = '''  +1.14%  '

# Row 4
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '''  -0.20%  '

# Row 5
$ws.Range('D5').Value = '''521.64'
$ws.Range('E5').Value = '''  -0.54%  '

# Row 6
$ws.Range('D6').Value = '''149.24'
$ws.Range('E6').Value = '''  -1.10%  '

# Row 7
$ws.Range('D7').Value = '''0.993'
$ws.Range('E7').Value = '''  -0.54%  '

# Row 8
$ws.Range('D8').Value = '''0.574'
$ws.Range('E8').Value = '''  -0.24%  '

# Row 9
$ws.Range('D9').Value = '''2.681.86'
$ws.Range('E9').Value = '''  +1.52%  '

# Row 10
$ws.Range('D10').Value = '''6.56'
$ws.Range('E10').Value = '''  +2.54%  '

# Row 11
$ws.Range('D11').Value = '''0.107'
$ws.Range('E11').Value = '''  +0.28%  '

# Row 12
$ws.Range('D12').Value = '''0.342'
$ws.Range('E12').Value = '''  -0.95%  '

# Row 13
$ws.Range('E13').Value = '''  -1.30%  '

# Row 14
$ws.Range('D14').Value = '''3.110.49'
$ws.Range('E14').Value = '''  +0.49%  '

# Row 15
$ws.Range('D15').Value = '''59.633.31'
$ws.Range('E15').Value = '''  -1.47%  '

# Row 16
$ws.Range('D16').Value = '''21.56'
$ws.Range('E16').Value = '''  +0.12%  '

# Row 17
$ws.Range('D17').Value = '''0.0000140'
$ws.Range('E17').Value = '''  +0.31%  '

# Row 18
$ws.Range('D18').Value = '''2.672.79'
$ws.Range('E18').Value = '''  +1.25%  '

# Row 19
$ws.Range('D19').Value = '''4.63'
$ws.Range('E19').Value = '''  -0.54%  '

# Row 20
$ws.Range('D20').Value = '''349.42'
$ws.Range('E20').Value = '''  +0.80%  '

# Row 21
$ws.Range('D21').Value = '''10.63'
$ws.Range('E21').Value = '''  +1.20%  '

# Row 22
$ws.Range('D22').Value = '''6.24'
$ws.Range('E22').Value = '''  +0.95%  '

# Row 23
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '''  +0.51%  '

# Row 24
$ws.Range('D24').Value = '''61.53'
$ws.Range('E24').Value = '''  +0.94%  '

# Row 25
$ws.Range('D25').Value = '''0.430'
$ws.Range('E25').Value = '''  +1.80%  '

# Row 26
$ws.Range('D26').Value = '''2.755.07'
$ws.Range('E26').Value = '''  +0.10%  '

# Row 27
$ws.Range('B27').Value = '''Kaspa'
$ws.Range('C27').Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '''0.162'
$ws.Range('E27').Value = '''  -1.21%  '

# Row 28
$ws.Range('B28').Value = '''Binance-PegBSC-USD'
$ws.Range('C28').Value = '''https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '''0.990'
$ws.Range('E28').Value = '''  -0.78%  '

# Row 29
$ws.Range('D29').Value = '''0.0₃0839'
$ws.Range('E29').Value = '''  +0.95%  '

# Row 30
$ws.Range('D30').Value = '''7.21'
$ws.Range('E30').Value = '''  +1.07%  '

# Row 31
$ws.Range('B31').Value = '''Aptos'
$ws.Range('C31').Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').Value = '''6.63'
$ws.Range('E31').Value = '''  +10.21%  '

# Row 32
$ws.Range('B32').Value = '''USDe'
$ws.Range('C32').Value = '''https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D32').Value = '''0.996'
$ws.Range('E32').Value = '''  -0.42%  '

# Row 33
$ws.Range('B33').Value = '''PancakeSwap'
$ws.Range('C33').Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '''1.59'
$ws.Range('E33').Value = '''  -0.48%  '

# Row 34
$ws.Range('B34').Value = '''EthereumClassic'
$ws.Range('C34').Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '''19.08'
$ws.Range('E34').Value = '''  -0.01%  '

# Row 35
$ws.Range('B35').Value = '''SuiNetwork'
$ws.Range('C35').Value = '''https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D35').Value = '''1.07'
$ws.Range('E35').Value = '''  +18.92%  '

# Row 36
$ws.Range('B36').Value = '''Monero'
$ws.Range('C36').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '''149.62'
$ws.Range('E36').Value = '''  -0.69%  '

# Row 37
$ws.Range('D37').Value = '''4.08'
$ws.Range('E37').Value = '''  +1.80%  '

# Row 38
$ws.Range('D38').Value = '''1.17'
$ws.Range('E38').Value = '''  +0.67%  '

# Row 39
$ws.Range('D39').Value = '''0.880'
$ws.Range('E39').Value = '''  +0.24%  '

# Row 40
$ws.Range('D40').Value = '''36.59'
$ws.Range('E40').Value = '''  -0.25%  '

# Row 41
$ws.Range('D41').Value = '''3.74'
$ws.Range('E41').Value = '''  +1.67%  '

# Row 42
$ws.Range('D42').Value = '''1.45'
$ws.Range('E42').Value = '''  -0.50%  '

# Row 43
$ws.Range('D43').Value = '''290.19'
$ws.Range('E43').Value = '''  -1.10%  '

# Row 44
$ws.Range('D44').Value = '''0.629'
$ws.Range('E44').Value = '''  -0.63%  '

# Row 45
$ws.Range('D45').Value = '''0.100'
$ws.Range('E45').Value = '''  -0.67%  '

# Row 46
$ws.Range('D46').Value = '''0.990'
$ws.Range('E46').Value = '''  -0.77%  '

# Row 47
$ws.Range('D47').Value = '''19.74'
$ws.Range('E47').Value = '''  -0.87%  '

# Row 48
$ws.Range('D48').Value = '''0.0550'
$ws.Range('E48').Value = '''  -0.47%  '

# Row 49
$ws.Range('B49').Value = '''RenderToken'
$ws.Range('C49').Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''4.83'
$ws.Range('E49').Value = '''  +1.02%  '

# Row 50
$ws.Range('B50').Value = '''VeChain'
$ws.Range('C50').Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '''0.0235'
$ws.Range('E50').Value = '''  -0.59%  '

# Row 51
$ws.Range('E51').Value = '''  -1.39%  '
